$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (prices / volume % changes), matching the
# upstream data refresh. A few coins (WrappedEther/ShibaInu, PEPE/USDe,
# Bittensor/FirstDigitalUSD) swapped rank position, so those rows have
# their Coin/Link/Price/Volume cells fully rewritten.

$ws.Range("D2").Value = "59.467.10"
$ws.Range("E2").Value = "  +2.76%  "
$ws.Range("D3").Value = "2.597.54"
$ws.Range("E3").Value = "  +1.17%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "536.46"
$ws.Range("E5").Value = "  +4.09%  "
$ws.Range("D6").Value = "140.94"
$ws.Range("E6").Value = "  +1.77%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "0.567"
$ws.Range("E8").Value = "  +1.55%  "
$ws.Range("D9").Value = "2.611.98"
$ws.Range("E9").Value = "  +1.21%  "
$ws.Range("D10").Value = "6.47"
$ws.Range("E10").Value = "  +0.31%  "
$ws.Range("E11").Value = "  +3.44%  "
$ws.Range("D12").Value = "0.335"
$ws.Range("E12").Value = "  +3.24%  "
$ws.Range("E13").Value = "  +3.05%  "
$ws.Range("D14").Value = "3.059.78"
$ws.Range("E14").Value = "  +1.18%  "
$ws.Range("D15").Value = "59.358.80"
$ws.Range("E15").Value = "  +2.58%  "
$ws.Range("D16").Value = "20.54"
$ws.Range("E16").Value = "  +2.54%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.608.00"
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.0000134"
$ws.Range("E18").Value = "  +2.14%  "
$ws.Range("D19").Value = "346.34"
$ws.Range("E19").Value = "  +4.16%  "
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("D21").Value = "10.15"
$ws.Range("E21").Value = "  +1.07%  "
$ws.Range("E22").Value = "  +0.62%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "67.16"
$ws.Range("E24").Value = "  +2.31%  "
$ws.Range("D25").Value = "0.169"
$ws.Range("E25").Value = "  +2.55%  "
$ws.Range("D26").Value = "0.407"
$ws.Range("E26").Value = "  +2.25%  "
$ws.Range("D27").Value = "0.997"
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").Value = "7.22"
$ws.Range("E28").Value = "  +4.77%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0739"
$ws.Range("E29").Value = "  +4.50%  "
$ws.Range("B30").Value = "USDe"
$ws.Range("C30").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").Value = "1.63"
$ws.Range("E31").Value = "  +5.36%  "
$ws.Range("D32").Value = "5.82"
$ws.Range("E32").Value = "  -1.14%  "
$ws.Range("D33").Value = "18.84"
$ws.Range("E33").Value = "  +1.45%  "
$ws.Range("D34").Value = "149.10"
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("D35").Value = "4.00"
$ws.Range("E35").Value = "  +3.05%  "
$ws.Range("E36").Value = "  +1.59%  "
$ws.Range("D37").Value = "36.95"
$ws.Range("E37").Value = "  +2.29%  "
$ws.Range("E38").Value = "  +4.98%  "
$ws.Range("D39").Value = "0.848"
$ws.Range("E39").Value = "  +2.67%  "
$ws.Range("D40").Value = "0.835"
$ws.Range("E40").Value = "  +0.92%  "
$ws.Range("D41").Value = "3.54"
$ws.Range("E41").Value = "  +1.48%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "277.89"
$ws.Range("E42").Value = "  +1.84%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("D44").Value = "0.599"
$ws.Range("E44").Value = "  +1.65%  "
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("E46").Value = "  +2.53%  "
$ws.Range("E47").Value = "  +1.69%  "
$ws.Range("D48").Value = "4.71"
$ws.Range("E48").Value = "  +5.06%  "
$ws.Range("D49").Value = "1.945.57"
$ws.Range("E49").Value = "  -1.50%  "
$ws.Range("D50").Value = "0.0223"
$ws.Range("E50").Value = "  +2.67%  "
$ws.Range("D51").Value = "18.34"
$ws.Range("E51").Value = "  +3.15%  "
